$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E12").Value = "['Normal']"

$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E25").Value = "['Normal']"

$ws.Range("D35").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E35").Value = "['Normal', 'HardwareFault']"

$ws.Range("D38").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['Normal', 'HardwareFault']"

$ws.Range("D54").Value = "[1, 0, 0, 0, 0, 1, 0]"
$ws.Range("E54").Value = "['Normal', 'CommunicationIssue']"

$ws.Range("D58").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['Normal', 'ParamViolation']"

$ws.Range("D67").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E67").Value = "['Normal', 'HardwareFault']"

$ws.Range("D69").Value = "[1, 1, 0, 0, 0, 1, 0]"
$ws.Range("E69").Value = "['Normal', 'SurroundingEnvironment', 'CommunicationIssue']"

$ws.Range("D70").Value = "[1, 1, 0, 0, 0, 1, 0]"
$ws.Range("E70").Value = "['Normal', 'SurroundingEnvironment', 'CommunicationIssue']"

$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

$ws.Range("D116").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E116").Value = "['Normal', 'SoftwareFault']"
